$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate text in cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")

$oldText = $cellA1.Value()
$newText = $oldText.Replace("1000 Bs = 12.59 = 51057.53 pesos", "1000 Bs = 12.57 = 51049.39 pesos")
$newText = $newText.Replace("51057.53 pesos = 12.54 = 970.45 Bs", "51049.39 pesos = 12.52 = 967.28 Bs")
$cellA1.Value = $newText

# --- Sheet "tasas": update the rate numbers ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 79.56999999999999
$wsTasas.Range("O10").Value = 4062
$wsTasas.Range("N12").Value = 4078
$wsTasas.Range("O12").Value = 77.27
